# Split the "Сделки" (Trades) sheet into two separate sheets:
#   - "Сделки с ЦБ"  (trades with securities / stocks)
#   - "Сделки с ПФИ" (trades with derivative financial instruments)
#
# The original "Сделки" sheet contained 4 instruments (2 rows each):
#   rows 10-11 : AAL   210115C00030000   (option  -> derivative)
#   rows 12-13 : VLO   200724P00064000   (option  -> derivative)
#   rows 14-15 : VLO                     (stock   -> security)
#   rows 16-17 : FB                      (stock   -> security)
#   row  18    : totals row
#
# Securities (VLO, FB) stay on the renamed "Сделки с ЦБ" sheet.
# Derivatives (the two options) move to the new "Сделки с ПФИ" sheet.

$wb = $excel.ActiveWorkbook

$trades = $wb.Worksheets.Item("Сделки")

# Create the new derivatives sheet as a full copy of the trades sheet
# (this preserves column widths, styles, merged cells and formulas),
# positioned immediately after the original sheet.
$trades.Copy($null, $trades)
$derivatives = $wb.Worksheets.Item(3)

# --- Rename sheets -------------------------------------------------
$trades.Name = "Сделки с ЦБ"
$derivatives.Name = "Сделки с ПФИ"

# --- "Сделки с ЦБ": keep only the securities (old rows 14-17) ------
# Deleting the first two instruments (rows 10-13, the options) shifts
# the remaining securities rows (14-17) up to rows 10-13, and the
# totals row (18) up to row 14; SUM() formulas auto-adjust.
$trades.Range("A10:P13").EntireRow.Delete()

# --- "Сделки с ПФИ": keep only the derivatives (old rows 10-13) ----
# Deleting the securities rows (14-17) removes the stock trades,
# leaving the two option instruments in rows 10-13, with the totals
# row shifting from 18 up to row 14.
$derivatives.Range("A14:P17").EntireRow.Delete()

# Update the report title and income/expense column headers on the
# derivatives sheet to the PFI-specific wording / tax codes.
$derivatives.Range("A1").Value = "Отчет по сделкам с производными финансовыми инструментами, завершённым в отчетном периоде"
$derivatives.Range("M8").Value = "Доход, RUB (код 1532)"
$derivatives.Range("N8").Value = "Расход, RUB (код 206)"

Write-Output "Sheets after split:"
foreach ($ws in $wb.Worksheets) {
    Write-Output $ws.Name
}
